$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "property" row (row 5) for the test_enum field, mirroring rows 1-4
$ws.Range("A5").Value = "#p"
$ws.Range("B5").Value = "test_enum"
$ws.Range("C5").Value = "enum_val_1"

# Extend the Data Dict header/data rows (row 7/8) with the new enum column (F)
$ws.Range("F7").Value = "test_enum"
$ws.Range("F8").Value = "enum_val_1"

# Update selection to match the new active cell
$ws.Range("C5").Select()
